# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (F column) counts to both the "展览"
# and "全部类型" sheets, which carry duplicate listings of the same
# convention data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 83
$ws.Range("F5").Value = 2743
$ws.Range("F9").Value = 1456
$ws.Range("F13").Value = 1223
$ws.Range("F15").Value = 374
$ws.Range("F22").Value = 2676
$ws.Range("F24").Value = 312

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 83
$ws.Range("F5").Value = 2743
$ws.Range("F9").Value = 1456
$ws.Range("F13").Value = 1223
$ws.Range("F15").Value = 374
$ws.Range("F20").Value = 77
$ws.Range("F22").Value = 2676
$ws.Range("F24").Value = 312
